$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    for ($col = 2; $col -le 29; $col++) {
        $v1 = $ws.Cells.Item($r1, $col).Value2
        $v2 = $ws.Cells.Item($r2, $col).Value2
        $ws.Cells.Item($r1, $col).Value = $v2
        $ws.Cells.Item($r2, $col).Value = $v1
    }
}

# Rows whose (non-id) data was re-ordered in the source feed
Swap-Rows 27 28
Swap-Rows 107 108
Swap-Rows 128 129
Swap-Rows 143 145

# Four brand-new fixtures are inserted before the old last row, which is
# pushed down to row 208 (and gains its previously-missing H/I/J/AB/AC values).
# Seed formatting (column A bordered/centered style, column E date format) by
# copying an existing fully-formatted row, then overwrite every value.
$ws.Range("A27:AC27").Copy()
$ws.Range("A204:AC208").PasteSpecial(-4122)

# Row 204
$ws.Cells.Item(204,1).Value = 202
$ws.Cells.Item(204,2).Value = 8010638
$ws.Cells.Item(204,3).Value = "Bolivia Primera División"
$ws.Cells.Item(204,4).Value = "Bolivia Apertura"
$ws.Cells.Item(204,5).Value = 45381.77083333334
$ws.Cells.Item(204,6).Value = "Guabira"
$ws.Cells.Item(204,7).Value = "Universitario De Vinto"
$ws.Cells.Item(204,8).Value = 3
$ws.Cells.Item(204,9).Value = 1
$ws.Cells.Item(204,10).Value = "H"
$ws.Cells.Item(204,11).Value = 2.2
$ws.Cells.Item(204,12).Value = 3.3
$ws.Cells.Item(204,13).Value = 3
$ws.Cells.Item(204,14).Value = 1.85
$ws.Cells.Item(204,15).Value = 3.5
$ws.Cells.Item(204,16).Value = 4.5
$ws.Cells.Item(204,17).Value = -0.5
$ws.Cells.Item(204,18).Value = 1.825
$ws.Cells.Item(204,19).Value = 1.975
$ws.Cells.Item(204,20).Value = 2.25
$ws.Cells.Item(204,21).Value = 1.825
$ws.Cells.Item(204,22).Value = 1.975
$ws.Cells.Item(204,23).Value = 0.8500000000000001
$ws.Cells.Item(204,24).Value = -1
$ws.Cells.Item(204,25).Value = -1
$ws.Cells.Item(204,26).Value = 0.825
$ws.Cells.Item(204,27).Value = -1
$ws.Cells.Item(204,28).Value = 0.825
$ws.Cells.Item(204,29).Value = -1

# Row 205
$ws.Cells.Item(205,1).Value = 203
$ws.Cells.Item(205,2).Value = 8010639
$ws.Cells.Item(205,3).Value = "Bolivia Primera División"
$ws.Cells.Item(205,4).Value = "Bolivia Apertura"
$ws.Cells.Item(205,5).Value = 45381.875
$ws.Cells.Item(205,6).Value = "Nacional Potosi"
$ws.Cells.Item(205,7).Value = "Royal Pari FC"
$ws.Cells.Item(205,8).Value = 2
$ws.Cells.Item(205,9).Value = 0
$ws.Cells.Item(205,10).Value = "H"
$ws.Cells.Item(205,11).Value = 1.363
$ws.Cells.Item(205,12).Value = 4.333
$ws.Cells.Item(205,13).Value = 7
$ws.Cells.Item(205,14).Value = 1.333
$ws.Cells.Item(205,15).Value = 5
$ws.Cells.Item(205,16).Value = 9.5
$ws.Cells.Item(205,17).Value = -1.5
$ws.Cells.Item(205,18).Value = 1.875
$ws.Cells.Item(205,19).Value = 1.925
$ws.Cells.Item(205,20).Value = 3.25
$ws.Cells.Item(205,21).Value = 1.9
$ws.Cells.Item(205,22).Value = 1.9
$ws.Cells.Item(205,23).Value = 0.333
$ws.Cells.Item(205,24).Value = -1
$ws.Cells.Item(205,25).Value = -1
$ws.Cells.Item(205,26).Value = 0.875
$ws.Cells.Item(205,27).Value = -1
$ws.Cells.Item(205,28).Value = -1
$ws.Cells.Item(205,29).Value = 0.8999999999999999

# Row 206
$ws.Cells.Item(206,1).Value = 204
$ws.Cells.Item(206,2).Value = 8010642
$ws.Cells.Item(206,3).Value = "Bolivia Primera División"
$ws.Cells.Item(206,4).Value = "Bolivia Apertura"
$ws.Cells.Item(206,5).Value = 45382.66666666666
$ws.Cells.Item(206,6).Value = "Always Ready"
$ws.Cells.Item(206,7).Value = "Independiente Petrolero"
$ws.Cells.Item(206,8).Value = 0
$ws.Cells.Item(206,9).Value = 1
$ws.Cells.Item(206,10).Value = "A"
$ws.Cells.Item(206,11).Value = 1.4
$ws.Cells.Item(206,12).Value = 4
$ws.Cells.Item(206,13).Value = 7
$ws.Cells.Item(206,14).Value = 1.3
$ws.Cells.Item(206,15).Value = 5.25
$ws.Cells.Item(206,16).Value = 10
$ws.Cells.Item(206,17).Value = -1.5
$ws.Cells.Item(206,18).Value = 1.875
$ws.Cells.Item(206,19).Value = 1.925
$ws.Cells.Item(206,20).Value = 3
$ws.Cells.Item(206,21).Value = 1.975
$ws.Cells.Item(206,22).Value = 1.825
$ws.Cells.Item(206,23).Value = -1
$ws.Cells.Item(206,24).Value = -1
$ws.Cells.Item(206,25).Value = 9
$ws.Cells.Item(206,26).Value = -1
$ws.Cells.Item(206,27).Value = 0.925
$ws.Cells.Item(206,28).Value = -1
$ws.Cells.Item(206,29).Value = 0.825

# Row 207
$ws.Cells.Item(207,1).Value = 205
$ws.Cells.Item(207,2).Value = 8010640
$ws.Cells.Item(207,3).Value = "Bolivia Primera División"
$ws.Cells.Item(207,4).Value = "Bolivia Apertura"
$ws.Cells.Item(207,5).Value = 45382.77083333334
$ws.Cells.Item(207,6).Value = "Bolivar"
$ws.Cells.Item(207,7).Value = "Oriente Petrolero"
$ws.Cells.Item(207,8).Value = 4
$ws.Cells.Item(207,9).Value = 2
$ws.Cells.Item(207,10).Value = "H"
$ws.Cells.Item(207,11).Value = 1.222
$ws.Cells.Item(207,12).Value = 5.75
$ws.Cells.Item(207,13).Value = 9
$ws.Cells.Item(207,14).Value = 1.166
$ws.Cells.Item(207,15).Value = 7
$ws.Cells.Item(207,16).Value = 15
$ws.Cells.Item(207,17).Value = -2.25
$ws.Cells.Item(207,18).Value = 2
$ws.Cells.Item(207,19).Value = 1.8
$ws.Cells.Item(207,20).Value = 3.25
$ws.Cells.Item(207,21).Value = 1.975
$ws.Cells.Item(207,22).Value = 1.825
$ws.Cells.Item(207,23).Value = 0.1659999999999999
$ws.Cells.Item(207,24).Value = -1
$ws.Cells.Item(207,25).Value = -1
$ws.Cells.Item(207,26).Value = -0.5
$ws.Cells.Item(207,27).Value = 0.4
$ws.Cells.Item(207,28).Value = 0.9750000000000001
$ws.Cells.Item(207,29).Value = -1

# Row 208 (was row 204 before the edit; now fully populated)
$ws.Cells.Item(208,1).Value = 206
$ws.Cells.Item(208,2).Value = 8011587
$ws.Cells.Item(208,3).Value = "Bolivia Primera División"
$ws.Cells.Item(208,4).Value = "Bolivia Apertura"
$ws.Cells.Item(208,5).Value = 45382.85416666666
$ws.Cells.Item(208,6).Value = "Club Aurora"
$ws.Cells.Item(208,7).Value = "Blooming"
$ws.Cells.Item(208,8).Value = 4
$ws.Cells.Item(208,9).Value = 2
$ws.Cells.Item(208,10).Value = "H"
$ws.Cells.Item(208,11).Value = 1.533
$ws.Cells.Item(208,12).Value = 3.75
$ws.Cells.Item(208,13).Value = 5.5
$ws.Cells.Item(208,14).Value = 1.615
$ws.Cells.Item(208,15).Value = 3.6
$ws.Cells.Item(208,16).Value = 6
$ws.Cells.Item(208,17).Value = -1
$ws.Cells.Item(208,18).Value = 2.025
$ws.Cells.Item(208,19).Value = 1.775
$ws.Cells.Item(208,20).Value = 2.5
$ws.Cells.Item(208,21).Value = 1.825
$ws.Cells.Item(208,22).Value = 1.975
$ws.Cells.Item(208,23).Value = 0.615
$ws.Cells.Item(208,24).Value = -1
$ws.Cells.Item(208,25).Value = -1
$ws.Cells.Item(208,26).Value = 1.025
$ws.Cells.Item(208,27).Value = -1
$ws.Cells.Item(208,28).Value = 0.825
$ws.Cells.Item(208,29).Value = -1
